# "Increase in download timeout"
#
# Adds a new DownloadTimeOut asset ("DownloadTimeOut" / "DownloadTimeOut_GPOAvendra")
# on the Assets sheet and records "Prod" as the OrchestratorAssetFolder for every
# existing asset row (column C), matching the pattern already used for every other
# asset in that sheet.

$wb = $excel.ActiveWorkbook

$wsAssets    = $wb.Worksheets.Item("Assets")
$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")

# Column C ("OrchestratorAssetFolder") is set to "Prod" for every existing asset
# row (2-13) as well as the brand-new row (14).
$wsAssets.Range("C2:C13").Value = "Prod"

# New asset: DownloadTimeOut / DownloadTimeOut_GPOAvendra, also in the "Prod" folder.
$wsAssets.Range("A14").Value = "DownloadTimeOut"
$wsAssets.Range("B14").Value = "DownloadTimeOut_GPOAvendra"
$wsAssets.Range("C14").Value = "Prod"

# Restore each sheet's on-screen selection (cursor position) as left after the edit.
# Select "Settings" last so it remains the active/visible tab when saved.
[void]$wsAssets.Range("C2:C14").Select()
[void]$wsConstants.Range("B37").Select()
[void]$wsSettings.Range("B8").Select()
